# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.795.31"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.917.60"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'356.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "'109.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.01%  "
$ws.Range("D7").Value = "'0.565"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "'39.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").Value = "'0.0870"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "'19.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").Value = "'7.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "3.373.84"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "2.925.54"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'0.980"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "51.790.16"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "'3.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").Value = "'7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "'13.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").Value = "'70.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").Value = "'268.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").Value = "'2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'0.185"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.31%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").Value = "'7.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +17.89%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'26.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +11.20%  "
$ws.Range("D31").Value = "'10.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "'37.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").Value = "'6.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "'2.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.67%  "
$ws.Range("D35").Value = "'52.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").Value = "'0.0440"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'3.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").Value = "'18.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("D41").Value = "'2.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.45%  "
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("D43").Value = "'22.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("D44").Value = "'118.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'3.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'2.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.04%  "
$ws.Range("D48").Value = "2.124.69"
$ws.Range("E48").Value = "  -3.51%  "
$ws.Range("E49").Value = "  -5.14%  "
$ws.Range("D50").Value = "'0.0337"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "'9.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.70%  "
